# Apply the FolderNames.xlsx edits:
#  - Fix a typo'd duplicate-ish entry:
#      "Rogue River-Siskiyou National Forests National Forest"
#      -> "Rogue River-Siskiyou National Forest"
#  - Remove the erroneous duplicate row:
#      "Siskiyou Mountains Ranger District Ranger District"
# Removing that row shifts everything below it up by one, which also
# removes the now-unused last row (previously "Woodruff Allotment" at A17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Rogue River-Siskiyou National Forests National Forest" text
$ws.Range("A11").Value = "Rogue River-Siskiyou National Forest"

# Remove the duplicate "Siskiyou Mountains Ranger District Ranger District" row entirely
$ws.Rows.Item(14).Delete()
